$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RateCompare")
$ws.Activate()

# Row 7
$ws.Range("A7").Value = "2112_GP10glc"
$ws.Range("B7").Value = "0-25"
$ws.Range("C7").Value = 0.17
$ws.Range("D7").Value = 0.01
$ws.Range("E7").Value = 0.53
$ws.Range("F7").Value = 0.04
$ws.Range("G7").Value = "EX_glc__D_e"

# Row 8
$ws.Range("A8").Value = "2112_GP20glc_Rep1"
$ws.Range("B8").Value = "0-40"
$ws.Range("C8").Value = 0.16
$ws.Range("D8").Value = 0.01
$ws.Range("E8").Value = 0.75
$ws.Range("F8").Value = 0.08
$ws.Range("G8").Value = "EX_glc__D_e"

# Row 9
$ws.Range("A9").Value = "2112_GP20glc_Rep2"
$ws.Range("B9").Value = "0-40"
$ws.Range("C9").Value = 0.16
$ws.Range("D9").Value = 0.01
$ws.Range("E9").Value = 0.76
$ws.Range("F9").Value = 0.07
$ws.Range("G9").Value = "EX_glc__D_e"

# Row 10
$ws.Range("A10").Value = "2112_GP20glc_Rep3"
$ws.Range("B10").Value = "0-40"
$ws.Range("C10").Value = 0.16
$ws.Range("D10").Value = 0.01
$ws.Range("E10").Value = 0.79
$ws.Range("F10").Value = 0.07

# Row 11
$ws.Range("A11").Value = "2112_GP50glc"
$ws.Range("B11").Value = "0-50"
$ws.Range("C11").Value = 0.12
$ws.Range("D11").Value = 0.02
$ws.Range("E11").Value = 0.3
$ws.Range("F11").Value = 0.19
$ws.Range("G11").Value = "EX_glc__D_e"

# Update selection to match the author's final cursor position
$ws.Range("F9").Select() | Out-Null
